$p = $ppt.ActivePresentation

# RDMPDEV-1336 Deleted old inaccurate documentation files
# Remove the three trailing slides (old inaccurate pipeline diagram
# documentation), keeping only the first slide in the deck.
for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}
